$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from row 24 (A24/D24) down to row 25, matching existing layout
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D24").Copy()
$ws.Range("D25").PasteSpecial(-4122) # xlPasteFormats

# Add new row 25 data
$ws.Range("A25").Value = 45440
$ws.Range("B25").Value = 6
$ws.Range("D25").Value = "Imputation into random forest seems to work the best…"

# Match row height from diff (rows with wrapped multi-line text get an explicit height)
$ws.Rows.Item(25).RowHeight = 28.5
